$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 112171795
$ws.Range("B9").Value = 77515
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = 'Garnlav'
$ws.Range("G9").Value = 'Alectoria sarmentosa'
$ws.Range("H9").Value = '(Ach.) Ach.'
$ws.Range("P9").Value = 'Mångberget, Skellefteå, Vb'
$ws.Range("Q9").Value = 756378.1091670797
$ws.Range("R9").Value = 7212049.955989202
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Västerbotten'
$ws.Range("U9").Value = 'Skellefteå'
$ws.Range("V9").Value = 'Västerbotten'
$ws.Range("W9").Value = 'Skellefteå socken'
$ws.Range("Y9").Value = '''2023-09-11'
$ws.Range("Z9").Value = '00:00'
$ws.Range("AA9").Value = '''2023-09-11'
$ws.Range("AB9").Value = '00:00'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = 'Carl Jansson'
$ws.Range("AX9").Value = 'Carl Jansson'
$ws.Range("AY9").Value = 'Länsstyrelsens naturvärdesinventeringar i Västerbottens län'

# Row 10
$ws.Range("A10").Value = 112171785
$ws.Range("B10").Value = 78542
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 229748
$ws.Range("F10").Value = 'Gytterlav'
$ws.Range("G10").Value = 'Protopannaria pezizoides'
$ws.Range("H10").Value = '(Weber) P.M.Jørg. & S.Ekman'
$ws.Range("P10").Value = 'Mångberget, Skellefteå, Vb'
$ws.Range("Q10").Value = 756412.4227988988
$ws.Range("R10").Value = 7211953.63443999
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Västerbotten'
$ws.Range("U10").Value = 'Skellefteå'
$ws.Range("V10").Value = 'Västerbotten'
$ws.Range("W10").Value = 'Skellefteå socken'
$ws.Range("Y10").Value = '''2023-09-11'
$ws.Range("Z10").Value = '00:00'
$ws.Range("AA10").Value = '''2023-09-11'
$ws.Range("AB10").Value = '00:00'
$ws.Range("AC10").Value = 'på berg'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = 'Carl Jansson'
$ws.Range("AX10").Value = 'Carl Jansson'
$ws.Range("AY10").Value = 'Länsstyrelsens naturvärdesinventeringar i Västerbottens län'

# Row 11
$ws.Range("A11").Value = 112171812
$ws.Range("B11").Value = 78611
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 6463
$ws.Range("F11").Value = 'Bårdlav'
$ws.Range("G11").Value = 'Nephroma parile'
$ws.Range("H11").Value = '(Ach.) Ach.'
$ws.Range("P11").Value = 'Mångberget, Skellefteå, Vb'
$ws.Range("Q11").Value = 756485.3358962236
$ws.Range("R11").Value = 7212023.397891168
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Västerbotten'
$ws.Range("U11").Value = 'Skellefteå'
$ws.Range("V11").Value = 'Västerbotten'
$ws.Range("W11").Value = 'Skellefteå socken'
$ws.Range("Y11").Value = '''2023-09-11'
$ws.Range("Z11").Value = '00:00'
$ws.Range("AA11").Value = '''2023-09-11'
$ws.Range("AB11").Value = '00:00'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AJ11").Value = 'sälg'
$ws.Range("AK11").Value = 'Salix caprea'
$ws.Range("AO11").Value = 'Salix caprea'
$ws.Range("AW11").Value = 'Carl Jansson'
$ws.Range("AX11").Value = 'Carl Jansson'
$ws.Range("AY11").Value = 'Länsstyrelsens naturvärdesinventeringar i Västerbottens län'

# Row 12
$ws.Range("A12").Value = 112171810
$ws.Range("B12").Value = 90332
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 4769
$ws.Range("F12").Value = 'Svavelriska'
$ws.Range("G12").Value = 'Lactarius scrobiculatus'
$ws.Range("H12").Value = '(Scop.:Fr.) Fr.'
$ws.Range("P12").Value = 'Mångberget, Skellefteå, Vb'
$ws.Range("Q12").Value = 756485.5908420115
$ws.Range("R12").Value = 7212020.437326429
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 'Västerbotten'
$ws.Range("U12").Value = 'Skellefteå'
$ws.Range("V12").Value = 'Västerbotten'
$ws.Range("W12").Value = 'Skellefteå socken'
$ws.Range("Y12").Value = '''2023-09-11'
$ws.Range("Z12").Value = '00:00'
$ws.Range("AA12").Value = '''2023-09-11'
$ws.Range("AB12").Value = '00:00'
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AW12").Value = 'Carl Jansson'
$ws.Range("AX12").Value = 'Carl Jansson'
$ws.Range("AY12").Value = 'Länsstyrelsens naturvärdesinventeringar i Västerbottens län'

# Row 13
$ws.Range("A13").Value = 112171779
$ws.Range("B13").Value = 101703
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 222412
$ws.Range("F13").Value = 'Tibast'
$ws.Range("G13").Value = 'Daphne mezereum'
$ws.Range("H13").Value = 'L.'
$ws.Range("P13").Value = 'Mångberget, Skellefteå, Vb'
$ws.Range("Q13").Value = 756291.0224138719
$ws.Range("R13").Value = 7211892.055915679
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 'Västerbotten'
$ws.Range("U13").Value = 'Skellefteå'
$ws.Range("V13").Value = 'Västerbotten'
$ws.Range("W13").Value = 'Skellefteå socken'
$ws.Range("Y13").Value = '''2023-09-11'
$ws.Range("Z13").Value = '00:00'
$ws.Range("AA13").Value = '''2023-09-11'
$ws.Range("AB13").Value = '00:00'
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AW13").Value = 'Carl Jansson'
$ws.Range("AX13").Value = 'Carl Jansson'
$ws.Range("AY13").Value = 'Länsstyrelsens naturvärdesinventeringar i Västerbottens län'

# Row 14
$ws.Range("A14").Value = 112171814
$ws.Range("B14").Value = 89351
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 112
$ws.Range("F14").Value = 'Stjärntagging'
$ws.Range("G14").Value = 'Asterodon ferruginosus'
$ws.Range("H14").Value = 'Pat.'
$ws.Range("P14").Value = 'Mångberget, Skellefteå, Vb'
$ws.Range("Q14").Value = 756486.3469059409
$ws.Range("R14").Value = 7212041.380007128
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 'Västerbotten'
$ws.Range("U14").Value = 'Skellefteå'
$ws.Range("V14").Value = 'Västerbotten'
$ws.Range("W14").Value = 'Skellefteå socken'
$ws.Range("Y14").Value = '''2023-09-11'
$ws.Range("Z14").Value = '00:00'
$ws.Range("AA14").Value = '''2023-09-11'
$ws.Range("AB14").Value = '00:00'
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AW14").Value = 'Carl Jansson'
$ws.Range("AX14").Value = 'Carl Jansson'
$ws.Range("AY14").Value = 'Länsstyrelsens naturvärdesinventeringar i Västerbottens län'

